$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Hyundai"
$ws.Range("A5").Value = "KIA"

$ws.Range("A6").Select()
